# Update EPEX spot prices workbook with the latest day (10-sep) and
# append new rows to the Gaz / CO2 sheets for 2025-09-08.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Prix Spot" -- add a new date column (CK) for "10-sep"
# ---------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("CK1").Value = "10-sep"
$wsPrix.Range("CJ1").Copy()
$wsPrix.Range("CK1").PasteSpecial(-4122)  # xlPasteFormats -- copy the bold/centered header style

$prixValues = @(8.550000000000001, 4.99, 16.02, 10.68, 6.81, 3.68, 30.01, 15, 20.62, 42.16, 8, 0, 0, -0.01, -0.01, -0.01, -0.01, 4.9, 23.09, 51.17, 41.79, 16.38, 14, 14.75)

for ($i = 0; $i -lt $prixValues.Count; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 89).Value = $prixValues[$i]
}

# ---------------------------------------------------------------
# Sheet 2: "Gaz" -- append row 86 for 2025-09-08
# ---------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A86").NumberFormat = "@"
$wsGaz.Range("A86").Value = "2025-09-08"
$wsGaz.Range("B86").Value = 32.175

# ---------------------------------------------------------------
# Sheet 3: "CO2" -- append row 86 for 2025-09-08 (no price yet)
# ---------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A86").NumberFormat = "@"
$wsCo2.Range("A86").Value = "2025-09-08"
$wsCo2.Range("B86").Value = ""
